# "Fruta / hortaliza, semanal" - weekly update: insert a new reporting row
# for Cilantro @ Vega Modelo de Temuco at the top of the data block (new
# row 357), pushing all subsequent rows (old 357-391) down by one to
# become 358-392.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 357; everything below (old rows 357-391) is
# pushed down to 358-392, carrying its formatting (incl. the date style
# on column D) along with it.
$ws.Rows.Item(357).Insert()

# Populate the newly inserted row 357 with this week's reading.
$ws.Cells.Item(357, 1).Value = 10
$ws.Cells.Item(357, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(357, 3).Value = "La Araucanía"
$ws.Cells.Item(357, 4).Value = 44769
$ws.Cells.Item(357, 5).Value = 9
$ws.Cells.Item(357, 6).Value = 100112040
$ws.Cells.Item(357, 7).Value = "Cilantro"
$ws.Cells.Item(357, 8).Value = "Sin especificar"
$ws.Cells.Item(357, 9).Value = "Primera"
$ws.Cells.Item(357, 10).Value = 100
$ws.Cells.Item(357, 11).Value = 4300
$ws.Cells.Item(357, 12).Value = 4300
$ws.Cells.Item(357, 13).Value = 4300
$ws.Cells.Item(357, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(357, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(357, 16).Value = 2150
$ws.Cells.Item(357, 17).Value = 2
$ws.Cells.Item(357, 18).Value = "Hortaliza"
